$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the figure labels / titles in column F (rows 2-6) to the new
# wording, replacing the "@V3" phrasing with "at 180-day visit" and
# moving it onto its own line.
$ws.Range("F2").Value = "CT abnormalities`nat 180-day visit"
$ws.Range("F3").Value = "CT Severity Score 1-5`nat 180-day visit"
$ws.Range("F4").Value = "CT Severity Score >5`nat 180-day visit"
$ws.Range("F5").Value = "Symptoms`nat 180-day visit"
$ws.Range("F6").Value = "Lung function impairment`nat 180-day visit"

# Leave the selection where the editor ended up after making the change.
[void]$ws.Range("F7").Select()
